$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- Update the "time_taken" timestamps (column F) on the "data" sheet ---
$data.Range("F2").Value  = "2021-10-05 14:22:14.949258"
$data.Range("F3").Value  = "2021-10-05 14:22:14.949267"
$data.Range("F4").Value  = "2021-10-05 14:22:14.949270"
$data.Range("F5").Value  = "2021-10-05 14:22:14.949273"
$data.Range("F6").Value  = "2021-10-05 14:22:14.949276"
$data.Range("F7").Value  = "2021-10-05 14:22:14.949278"
$data.Range("F8").Value  = "2021-10-05 14:22:14.949281"
$data.Range("F9").Value  = "2021-10-05 14:22:14.949284"
$data.Range("F10").Value = "2021-10-05 14:22:14.949287"
$data.Range("F11").Value = "2021-10-05 14:22:14.949290"
$data.Range("F12").Value = "2021-10-05 14:22:14.949292"
$data.Range("F13").Value = "2021-10-05 14:22:14.949295"
$data.Range("F14").Value = "2021-10-05 14:22:14.949297"
$data.Range("F15").Value = "2021-10-05 14:22:14.949300"
$data.Range("F16").Value = "2021-10-05 14:22:14.949303"
$data.Range("F17").Value = "2021-10-05 14:22:14.949306"
$data.Range("F18").Value = "2021-10-05 14:22:14.949309"
$data.Range("F19").Value = "2021-10-05 14:22:14.949311"
$data.Range("F20").Value = "2021-10-05 14:22:14.949314"
$data.Range("F21").Value = "2021-10-05 14:22:14.949317"

# --- Add the new "metadata" worksheet right after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Reuse the header style (bold, centered, thin border) already defined in the
# workbook by copying the formatting from an existing styled header cell.
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Polycystic liver disease interim"
$meta.Range("C2").Value = 653

# data_version must be stored as text "1.23" (not the number 1.23)
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.23"
$meta.Range("D2").ClearFormats()

$meta.Range("E2").Value = "2021-03-15T15:30:01.022098Z"
$meta.Range("F2").Value = "2021-10-05 14:22:14.945827"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/653/?format=json"
